# Sprint1 Amend - Review phase in scrum board must have effort attributed
# to it in burndown (they count as part of the task).
#
# The scrum board review phase effort now gets split off into its own
# cell next to the task's existing day entry, so several existing day
# values shrink and the freed-up effort lands in an adjacent day column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Burndown Chart")

# Task "T4" row (row 8): Dia 6 (J) effort reduced, remainder now recorded
# under Dia 7 (K) for the review phase.
$ws.Range("J8").Value = 0.1
$ws.Range("K8").Value = 0.1

# Task "T5" row (row 9): Dia 7 (K) effort reduced, remainder now recorded
# under Dia 8 (L) for the review phase.
$ws.Range("K9").Value = 0.5
$ws.Range("L9").Value = 0.5

# task row (row 10): Dia 7 (K) effort reduced, remainder now recorded
# under Dia 11 (O) for the review phase.
$ws.Range("K10").Value = 0.2
$ws.Range("O10").Value = 0.3

# Task "T6" row (row 11): Dia 8 (L) effort reduced, remainder now recorded
# under Dia 9 (M) for the review phase.
$ws.Range("L11").Value = 0.7
$ws.Range("M11").Value = 0.3

# Update the view so the active cell / scroll position reflects where the
# edits were made.
$ws.Activate()
$ws.Range("O11").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.Zoom = 67
